$d = $word.ActiveDocument

$replacements = @(
    ,@("2024-05-19 Sunday", "2024-05-20 Monday")
    ,@("15×38=570", "34×38=1292")
    ,@("75×37=2775", "39×49=1911")
    ,@("66×32=2112", "74×88=6512")
    ,@("21×95=1995", "84×42=3528")
    ,@("14×61=854", "47×94=4418")
    ,@("41×42=1722", "78×53=4134")
    ,@("62×29=1798", "60×50=3000")
    ,@("62×60=3720", "19×39=741")
    ,@("83×84=6972", "38×61=2318")
    ,@("74×60=4440", "65×72=4680")
    ,@("57×19=1083", "82×84=6888")
    ,@("61×52=3172", "38×83=3154")
    ,@("69×20=1380", "29×38=1102")
    ,@("12×98=1176", "95×91=8645")
    ,@("55×45=2475", "74×35=2590")
    ,@("16×44=704", "15×96=1440")
    ,@("51×15=765", "96×86=8256")
    ,@("41×60=2460", "72×28=2016")
    ,@("86×87=7482", "48×98=4704")
    ,@("80×29=2320", "55×78=4290")
    ,@("83×28=2324", "31×65=2015")
    ,@("56×48=2688", "65×60=3900")
    ,@("70×99=6930", "80×75=6000")
    ,@("83×23=1909", "84×73=6132")
    ,@("71×22=1562", "22×42=924")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Output "NOT FOUND: $old"
    }
}

$d.Save()
